$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 (columns B-G) with new values
$ws.Range("B2").Value = 0.5224929946896146
$ws.Range("C2").Value = 1.898968271571221
$ws.Range("D2").Value = 15.67148254258704
$ws.Range("E2").Value = 3.958722337142003
$ws.Range("F2").Value = 3.967451592372942
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = 0.4961555363655573
$ws.Range("C3").Value = 1.850890984938442
$ws.Range("D3").Value = 15.31652666909897
$ws.Range("E3").Value = 3.913633435708941
$ws.Range("F3").Value = 3.925922114836308
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = 0.6562009383051112
$ws.Range("C4").Value = 1.755547849851764
$ws.Range("D4").Value = 14.91514394174393
$ws.Range("E4").Value = 3.862012939095871
$ws.Range("F4").Value = 3.849856378757894
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = 0.5233911015648834
$ws.Range("C5").Value = 1.855814880858909
$ws.Range("D5").Value = 15.64356210993208
$ws.Range("E5").Value = 3.955194320122853
$ws.Range("F5").Value = 3.966808199306541
$ws.Range("G5").Value = 43

# Add new rows 6-11 with labels Q4-Q9 and values, copying style from A5
$ws.Range("A6").Value = "Q4"
$ws.Range("B6").Value = 0.7469272014119552
$ws.Range("C6").Value = 1.819658143763748
$ws.Range("D6").Value = 15.50982749350374
$ws.Range("E6").Value = 3.938251832159003
$ws.Range("F6").Value = 3.913643877732932
$ws.Range("G6").Value = 42

$ws.Range("A7").Value = "Q5"
$ws.Range("B7").Value = 0.6021480811658927
$ws.Range("C7").Value = 1.875207305608213
$ws.Range("D7").Value = 16.24141422234426
$ws.Range("E7").Value = 4.030063798793297
$ws.Range("F7").Value = 4.034328036793712
$ws.Range("G7").Value = 41

$ws.Range("A8").Value = "Q6"
$ws.Range("B8").Value = 0.7822106663771212
$ws.Range("C8").Value = 1.836403048259028
$ws.Range("D8").Value = 16.36623743788394
$ws.Range("E8").Value = 4.045520663386103
$ws.Range("F8").Value = 4.019744080550121
$ws.Range("G8").Value = 40

$ws.Range("A9").Value = "Q7"
$ws.Range("B9").Value = 0.5940989821051248
$ws.Range("C9").Value = 1.932689778747934
$ws.Range("D9").Value = 16.9873845567875
$ws.Range("E9").Value = 4.121575494490851
$ws.Range("F9").Value = 4.131849360674751
$ws.Range("G9").Value = 39

$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.830600775663463
$ws.Range("C10").Value = 1.865149275164063
$ws.Range("D10").Value = 17.27800395215996
$ws.Range("E10").Value = 4.156681843990463
$ws.Range("F10").Value = 4.12752147190327
$ws.Range("G10").Value = 38

$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.5950944418594691
$ws.Range("C11").Value = 1.827659872411423
$ws.Range("D11").Value = 17.20024814321981
$ws.Range("E11").Value = 4.147318186879301
$ws.Range("F11").Value = 4.161016494713655
$ws.Range("G11").Value = 37

# Apply the same style as A2:A5 (bold, bordered, centered) to the new labels A6:A11
$ws.Range("A2").Copy()
$ws.Range("A6:A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
